# Update cached pricing/profit values on several Leve-profit worksheets.
# Source data refreshed by the scheduled market-data runner; only the
# cached numeric results (columns H-N) change, no formulas are involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 39999.5
$ws.Range("J3").Value = 39999.5
$ws.Range("L3").Value = 39999.5
$ws.Range("N3").Value = -40227.5

# Row 86
$ws.Range("H86").Value = 4167.423
$ws.Range("I86").Value = 1439.6842
$ws.Range("J86").Value = 11571.286
$ws.Range("K86").Value = 1439.6842
$ws.Range("L86").Value = 11571.286
$ws.Range("M86").Value = -316.6841999999999
$ws.Range("N86").Value = -13817.286

# Row 89
$ws.Range("H89").Value = 4167.423
$ws.Range("I89").Value = 1439.6842
$ws.Range("J89").Value = 11571.286
$ws.Range("K89").Value = 7198.420999999999
$ws.Range("L89").Value = 57856.43
$ws.Range("M89").Value = -1582.420999999999
$ws.Range("N89").Value = -69088.42999999999

# Row 92
$ws.Range("H92").Value = 2066.5557
$ws.Range("I92").Value = 3079.8235
$ws.Range("J92").Value = 344
$ws.Range("K92").Value = 3079.8235
$ws.Range("L92").Value = 344
$ws.Range("M92").Value = -1831.8235
$ws.Range("N92").Value = -2840

# Row 102
$ws.Range("H102").Value = 39999.5
$ws.Range("J102").Value = 39999.5
$ws.Range("L102").Value = 39999.5
$ws.Range("N102").Value = -46489.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 284405.6
$ws.Range("I61").Value = 201093.36
$ws.Range("J61").Value = 482768.1
$ws.Range("K61").Value = 201093.36
$ws.Range("L61").Value = 482768.1
$ws.Range("M61").Value = -200881.36
$ws.Range("N61").Value = -483192.1

# Row 74
$ws.Range("H74").Value = 9967162
$ws.Range("I74").Value = 7464428.5
$ws.Range("J74").Value = 15200150
$ws.Range("K74").Value = 7464428.5
$ws.Range("L74").Value = 15200150
$ws.Range("M74").Value = -7463554.5
$ws.Range("N74").Value = -15201898

# Row 77
$ws.Range("H77").Value = 9967162
$ws.Range("I77").Value = 7464428.5
$ws.Range("J77").Value = 15200150
$ws.Range("K77").Value = 37322142.5
$ws.Range("L77").Value = 76000750
$ws.Range("M77").Value = -37317774.5
$ws.Range("N77").Value = -76009486

# Row 132
$ws.Range("H132").Value = 26340.146
$ws.Range("I132").Value = 36977.535
$ws.Range("J132").Value = 3428.8462
$ws.Range("K132").Value = 110932.605
$ws.Range("L132").Value = 10286.5386
$ws.Range("M132").Value = -108402.605
$ws.Range("N132").Value = -15346.5386

# Row 136
$ws.Range("H136").Value = 284405.6
$ws.Range("I136").Value = 201093.36
$ws.Range("J136").Value = 482768.1
$ws.Range("K136").Value = 603280.08
$ws.Range("L136").Value = 1448304.3
$ws.Range("M136").Value = -600730.08
$ws.Range("N136").Value = -1453404.3

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 903.65
$ws.Range("I94").Value = 797.3333
$ws.Range("J94").Value = 1063.125
$ws.Range("K94").Value = 797.3333
$ws.Range("L94").Value = 1063.125
$ws.Range("M94").Value = -346.3333
$ws.Range("N94").Value = -1965.125

# Row 105
$ws.Range("H105").Value = 29461.234
$ws.Range("I105").Value = 48768
$ws.Range("J105").Value = 1880.1428
$ws.Range("K105").Value = 48768
$ws.Range("L105").Value = 1880.1428
$ws.Range("M105").Value = -47021
$ws.Range("N105").Value = -5374.1428

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1436.29
$ws.Range("I31").Value = 744.0612
$ws.Range("J31").Value = 2101.3726
$ws.Range("K31").Value = 744.0612
$ws.Range("L31").Value = 2101.3726
$ws.Range("M31").Value = -449.0612
$ws.Range("N31").Value = -2691.3726

# Row 34
$ws.Range("H34").Value = 1436.29
$ws.Range("I34").Value = 744.0612
$ws.Range("J34").Value = 2101.3726
$ws.Range("K34").Value = 744.0612
$ws.Range("L34").Value = 2101.3726
$ws.Range("M34").Value = -542.0612
$ws.Range("N34").Value = -2505.3726

# Row 86
$ws.Range("H86").Value = 2884.9429
$ws.Range("I86").Value = 2012.409
$ws.Range("J86").Value = 4361.5386
$ws.Range("K86").Value = 2012.409
$ws.Range("L86").Value = 4361.5386
$ws.Range("M86").Value = -889.4090000000001
$ws.Range("N86").Value = -6607.5386

# Row 89
$ws.Range("H89").Value = 2884.9429
$ws.Range("I89").Value = 2012.409
$ws.Range("J89").Value = 4361.5386
$ws.Range("K89").Value = 10062.045
$ws.Range("L89").Value = 21807.693
$ws.Range("M89").Value = -4446.045
$ws.Range("N89").Value = -33039.693

# Row 141
$ws.Range("H141").Value = 19584.285
$ws.Range("J141").Value = 19584.285
$ws.Range("L141").Value = 19584.285
$ws.Range("N141").Value = -29944.285

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1170.7858
$ws.Range("J68").Value = 1971.4762
$ws.Range("L68").Value = 5914.4286
$ws.Range("N68").Value = -7536.4286

# Row 71
$ws.Range("H71").Value = 1170.7858
$ws.Range("J71").Value = 1971.4762
$ws.Range("L71").Value = 17743.2858
$ws.Range("N71").Value = -25855.2858

# Row 131
$ws.Range("H131").Value = 1242.7451
$ws.Range("J131").Value = 1283.7693
$ws.Range("L131").Value = 3851.3079
$ws.Range("N131").Value = -13931.3079

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 3934633.5
$ws.Range("I14").Value = 4215571.5
$ws.Range("J14").Value = 1500
$ws.Range("K14").Value = 4215571.5
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = -4215403.5
$ws.Range("N14").Value = -1836

# Row 70
$ws.Range("H70").Value = 7952.1562
$ws.Range("I70").Value = 4463.95
$ws.Range("J70").Value = 13765.833
$ws.Range("K70").Value = 4463.95
$ws.Range("L70").Value = 13765.833
$ws.Range("M70").Value = -4193.95
$ws.Range("N70").Value = -14305.833

# Row 73
$ws.Range("H73").Value = 7952.1562
$ws.Range("I73").Value = 4463.95
$ws.Range("J73").Value = 13765.833
$ws.Range("K73").Value = 4463.95
$ws.Range("L73").Value = 13765.833
$ws.Range("M73").Value = -3527.95
$ws.Range("N73").Value = -15637.833

# Row 132
$ws.Range("H132").Value = 2361872
$ws.Range("I132").Value = 3790944
$ws.Range("J132").Value = 3903.4
$ws.Range("K132").Value = 11372832
$ws.Range("L132").Value = 11710.2
$ws.Range("M132").Value = -11370302
$ws.Range("N132").Value = -16770.2

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2682.2031
$ws.Range("I136").Value = 1319.2157
$ws.Range("J136").Value = 8029.3076
$ws.Range("K136").Value = 3957.6471
$ws.Range("L136").Value = 24087.9228
$ws.Range("M136").Value = -1407.6471
$ws.Range("N136").Value = -29187.9228

$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 19313.334
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 19313.334
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 19313.334
$ws.Range("N58").Value = -19929.334
$ws.Range("M58").ClearContents()
